$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New commit row: "talk virtual function instruction repaired" with 0.5 hours logged
$ws.Range("C22").Value = "talk virtual function instruction repaired"
$ws.Range("C22").Style = $ws.Range("C21").Style
$ws.Range("G22").Value = 0.5

# Move the active selection from E23 to D23 (matches the saved sheetView)
$ws.Range("D23").Select()
